$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The 4-cylinder group (previously a single data row, row 7) needs to
# be split into two rows (engine=0 and engine=1), just like the
# 6-cylinder group already is. This pushes every following row down
# by one and extends the table by one row.
# ------------------------------------------------------------------

# Insert a new row above the current row 8 ("6 cylinder" row). This
# shifts rows 8-11 down to 9-12 and correctly shifts/extends the
# existing merged cells (A8:A9 -> A9:A10, A11:G11 -> A12:G12).
$ws.Rows.Item(8).Insert()

# Copy the formatting of the (still unmodified) row 7 onto the newly
# inserted row 8 so the new row starts from the same base style as
# the rest of the data rows.
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H8").PasteSpecial(-4122)
$ws.Range("A7:H7").Copy()
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Write the new data values (mtcars data set) for rows 7-11.
# ------------------------------------------------------------------

# Row 7: 4 cylinder / engine 0
$ws.Cells.Item(7,1).Value = 4
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(7,4).Value = 91
$ws.Cells.Item(7,5).Value = $null
$ws.Cells.Item(7,6).Value = 2.14
$ws.Cells.Item(7,7).Value = $null

# Row 8: 4 cylinder / engine 1
$ws.Cells.Item(8,1).Value = $null
$ws.Cells.Item(8,2).Value = 1
$ws.Cells.Item(8,3).Value = 10
$ws.Cells.Item(8,4).Value = 81.8
$ws.Cells.Item(8,5).Value = 21.87235698318771
$ws.Cells.Item(8,6).Value = 2.3003
$ws.Cells.Item(8,7).Value = 0.5982073312080948

# Row 9: 6 cylinder / engine 0
$ws.Cells.Item(9,1).Value = 6
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 3
$ws.Cells.Item(9,4).Value = 131.6666666666667
$ws.Cells.Item(9,5).Value = 37.52776749732568
$ws.Cells.Item(9,6).Value = 2.755
$ws.Cells.Item(9,7).Value = 0.1281600561797629

# Row 10: 6 cylinder / engine 1
$ws.Cells.Item(10,1).Value = $null
$ws.Cells.Item(10,2).Value = 1
$ws.Cells.Item(10,3).Value = 4
$ws.Cells.Item(10,4).Value = 115.25
$ws.Cells.Item(10,5).Value = 9.178779875342908
$ws.Cells.Item(10,6).Value = 3.38875
$ws.Cells.Item(10,7).Value = 0.1162163929916946

# Row 11: 8 cylinder / engine 0 (single row, no split)
$ws.Cells.Item(11,1).Value = 8
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 14
$ws.Cells.Item(11,4).Value = 209.2142857142857
$ws.Cells.Item(11,5).Value = 50.97688551827051
$ws.Cells.Item(11,6).Value = 3.999214285714287
$ws.Cells.Item(11,7).Value = 0.7594047444769265

# ------------------------------------------------------------------
# Fix up the "N" column (cyl/engine group label) styling/merges:
# the A column uses a "top aligned" style on the first row of a
# merged 2-row group, and the plain style on the single-row (8 cyl)
# group - matching the layout already used for the 6-cylinder group.
# ------------------------------------------------------------------

# A7 (top of new A7:A8 merge) should look like A9 (top of A9:A10 merge)
$ws.Range("A9").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Merge A7:A8 together, for the newly split 4-cylinder group.
$ws.Range("A7:A8").Merge()

# A8 (bottom of the new merge) should look like A10/H-column plain style
# (the bottom cell of a merged N-column group carries no special
# border/alignment formatting). This must happen AFTER the merge,
# since merging re-applies the top-left cell's style across the
# whole merged range.
$ws.Range("H8").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Done"
